# Setzen des Status zu einem Flug (FlugHandler getFlugStatus()) & erste Dokumentation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row for the "getFlugstatus" documentation entry, right below the
#    "getAllFlugModels" row (old row 12), pushing everything from row 13 onward down by one.
$ws.Rows(13).Insert()

$ws.Cells.Item(13, 2).Value = "getFlugstatus"
$ws.Cells.Item(13, 3).Value = "String aktuellString, String flugstring"
$ws.Cells.Item(13, 4).Value = "String"
$ws.Cells.Item(13, 5).ClearFormats()
$ws.Cells.Item(13, 5).Value = """Der Status kann nicht gesetzt werden."""
$ws.Cells.Item(13, 6).Value = "-"
$ws.Cells.Item(13, 7).Value = "Ausgabe des Status eines Flugs zu einer gegebeben Uhrzeit"

# 2) Update the "assignFlugzeugToFlug" row (old row 16, now shifted to row 17):
#    new / extended error message and method description, with wrapped text and a taller row.
$ws.Cells.Item(17, 5).Value = " - ""Dieses Flugzeug wurde bereits einem Flug zugeordnet.""                              - ""Dieser Flugzeug wurde bereits diesem Flug zugeordnet."""
$ws.Cells.Item(17, 5).WrapText = $true

$ws.Cells.Item(17, 7).Value = "Methode, die ein Flugzeug einem Flug zuordnet. Zu beachten ist, dass ein Flugzeug nur genau einem Flug zugeordnet werden kann!"
$ws.Cells.Item(17, 7).WrapText = $true

$ws.Rows(17).RowHeight = 60

# 3) Widen column E a bit to better fit the new, longer text.
$ws.Columns("E").ColumnWidth = 37.8

# 4) Update the view so the selected cell is the new documentation row (G13).
$ws.Application.Goto($ws.Range("A7"))
$ws.Range("G13").Select()
